# api token and value code changes
# Rotates the /user/* and /project/* CRUD rows to new order/values and
# refreshes the shared test-data tokens (userID/actionID/projectID and the
# random name/email/description strings) used throughout the report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- /user section (rows 3-6 rotate: create -> row3, get -> row4, update -> row5, delete -> row6) ---
$ws.Range("G2").Value = '{"userID":"670"}'

$ws.Range("B3").Value = "POST"
$ws.Range("D3").Value = "/user/createUser"
$ws.Range("E3").Value = '{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"}'

$ws.Range("B4").Value = "GET"
$ws.Range("D4").Value = "/user/getUsers"
$ws.Range("G4").ClearContents()

$ws.Range("B5").Value = "PUT"
$ws.Range("D5").Value = "/user/updateUser"
$ws.Range("E5").Value = '{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"}'

$ws.Range("B6").Value = "DELETE"
$ws.Range("D6").Value = "/user/deleteUser"
$ws.Range("G6").Value = '{"userID":"670"}'

# --- /action section ---
$ws.Range("E8").Value = '{"updatedBy":{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"},"createdBy":{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"},"actionID":"476","description":"aKlOb","projectID":{"updatedBy":{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"},"createdBy":{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"},"description":"5rGdL","projectName":"0NwWb","projectID":"198"},"actionName":"hp0RM"}'

$ws.Range("E9").Value = '{"updatedBy":{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"},"createdBy":{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"},"actionID":"476","description":"aKlOb","projectID":{"updatedBy":{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"},"createdBy":{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"},"description":"5rGdL","projectName":"0NwWb","projectID":"198"},"actionName":"hp0RM"}'

$ws.Range("G10").Value = '{"actionID":"476"}'
$ws.Range("G11").Value = '{"actionID":"476"}'

# --- /project section ---
$ws.Range("E13").Value = '{"updatedBy":{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"},"createdBy":{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"},"description":"5rGdL","projectName":"0NwWb","projectID":"198"}'

# Row 14/15 swap: update -> row15, delete -> row14 (with refreshed values)
$ws.Range("B14").Value = "DELETE"
$ws.Range("D14").Value = "/project/deleteProject"
$ws.Range("E14").ClearContents()
$ws.Range("G14").Value = '{"projectID":"198"}'
$ws.Range("N14").ClearContents()

$ws.Range("B15").Value = "PUT"
$ws.Range("D15").Value = "/project/updateProject"
$ws.Range("E15").Value = '{"updatedBy":{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"},"createdBy":{"firstName":"RE8hP","lastName":"9k60a","emailID":"jESW3","userID":"670"},"description":"5rGdL","projectName":"0NwWb","projectID":"198"}'
$ws.Range("G15").ClearContents()
$ws.Range("N15").Value = "Project"

$ws.Range("F16").Value = '{"projectID":"198"}'
